# Insert a new NetCall entry ("SendServerInfo") as row 10 of the
# "netcalls" worksheet / "Table2" table, shifting every following row
# down by one (table grows from A1:E94 to A1:E95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at worksheet row 10 (row 9 of the table body,
# right after "SendWarnRequest" / before "SendUnbanRequest").
$ws.Rows.Item(10).Insert()

# Populate the new row with the new NetCall's data.
$ws.Cells.Item(10, 1).Value = "SendServerInfo"
$ws.Cells.Item(10, 2).Value = 1008
$ws.Cells.Item(10, 3).Value = "Data.NetCall"
$ws.Cells.Item(10, 4).Value = "FROM_CLIENT"
$ws.Cells.Item(10, 5).Value = "WarfareServerInfo info"

# Grow the table/list-object so it covers the newly inserted row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E95"))

# Match the author's final selection/view state.
$ws.Range("E10").Select() | Out-Null
